$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = 352
$ws.Range("D6").Value = 293
$ws.Range("E6").Value = 59
$ws.Range("F6").Value = 64.11378555798687
$ws.Range("G6").Value = 16.76136363636364
$ws.Range("H6").Value = 83.23863636363636
